$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---
# Copy formatting from the existing header (H1) so the new headers match
# the bold / centered / bordered style used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- New data columns I2:I32 and J2:J32 ---
$iValues = @(6,7,7,7,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,7,1,1,2,2,3,6,5,5,6,1,1)
$jValues = @(8,9,8,9,4,2,6,5,7,5,6,4,7,6,7,1,3,6,5,9,5,4,4,5,5,9,8,8,8,3,2)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
